$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2318.4546
$ws.Range("J32").Value = 2480.3
$ws.Range("L32").Value = 2480.3
$ws.Range("N32").Value = -3132.3
$ws.Range("H40").Value = 2732.6667
$ws.Range("I40").Value = 1466.3334
$ws.Range("J40").Value = 3999
$ws.Range("K40").Value = 1466.3334
$ws.Range("L40").Value = 3999
$ws.Range("M40").Value = -1291.3334
$ws.Range("N40").Value = -4349
$ws.Range("H43").Value = 9001.333000000001
$ws.Range("I43").Value = 9997.5
$ws.Range("K43").Value = 9997.5
$ws.Range("M43").Value = -9928.5
$ws.Range("H70").Value = 3058.0715
$ws.Range("I70").Value = 2939.25
$ws.Range("J70").Value = 3216.5
$ws.Range("K70").Value = 8817.75
$ws.Range("L70").Value = 9649.5
$ws.Range("M70").Value = -8547.75
$ws.Range("N70").Value = -10189.5
$ws.Range("H73").Value = 3058.0715
$ws.Range("I73").Value = 2939.25
$ws.Range("J73").Value = 3216.5
$ws.Range("K73").Value = 8817.75
$ws.Range("L73").Value = 9649.5
$ws.Range("M73").Value = -7881.75
$ws.Range("N73").Value = -11521.5
$ws.Range("H74").Value = 7035.25
$ws.Range("I74").Value = 6595.5
$ws.Range("J74").Value = 7475
$ws.Range("K74").Value = 6595.5
$ws.Range("L74").Value = 7475
$ws.Range("M74").Value = -5659.5
$ws.Range("N74").Value = -9347
$ws.Range("H77").Value = 7035.25
$ws.Range("I77").Value = 6595.5
$ws.Range("J77").Value = 7475
$ws.Range("K77").Value = 32977.5
$ws.Range("L77").Value = 37375
$ws.Range("M77").Value = -28297.5
$ws.Range("N77").Value = -46735
$ws.Range("H101").Value = 367.54544
$ws.Range("I101").Value = 238.5
$ws.Range("J101").Value = 711.6667
$ws.Range("K101").Value = 715.5
$ws.Range("L101").Value = 2135.0001
$ws.Range("M101").Value = 906.5
$ws.Range("N101").Value = -5379.0001
$ws.Range("H137").Value = 13133.567
$ws.Range("I137").Value = 2522.35
$ws.Range("K137").Value = 7567.049999999999
$ws.Range("M137").Value = -5017.049999999999
$ws.Range("H138").Value = 3328.7917
$ws.Range("I138").Value = 2412.3076
$ws.Range("J138").Value = 3669.2
$ws.Range("K138").Value = 7236.9228
$ws.Range("L138").Value = 11007.6
$ws.Range("M138").Value = -2096.9228
$ws.Range("N138").Value = -21287.6
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 898339.75
$ws.Range("I61").Value = 3670.6538
$ws.Range("K61").Value = 3670.6538
$ws.Range("M61").Value = -3458.6538
$ws.Range("H136").Value = 898339.75
$ws.Range("I136").Value = 3670.6538
$ws.Range("K136").Value = 11011.9614
$ws.Range("M136").Value = -8461.9614
$ws.Range("H141").Value = 57990
$ws.Range("J141").Value = 57990
$ws.Range("L141").Value = 57990
$ws.Range("N141").Value = -68350
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 85796330
$ws.Range("I132").Value = 4902.6665
$ws.Range("K132").Value = 14707.9995
$ws.Range("M132").Value = -12177.9995
$ws.Range("H140").Value = 126196.93
$ws.Range("J140").Value = 120672.68
$ws.Range("L140").Value = 120672.68
$ws.Range("N140").Value = -131032.68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2610.682
$ws.Range("I11").Value = 1675.1428
$ws.Range("J11").Value = 3047.2666
$ws.Range("K11").Value = 5025.428400000001
$ws.Range("L11").Value = 9141.799800000001
$ws.Range("M11").Value = -4885.428400000001
$ws.Range("N11").Value = -9421.799800000001
$ws.Range("H18").Value = 371
$ws.Range("I18").Value = 279.83334
$ws.Range("J18").Value = 644.5
$ws.Range("K18").Value = 839.5000200000001
$ws.Range("L18").Value = 1933.5
$ws.Range("M18").Value = -670.5000200000001
$ws.Range("N18").Value = -2271.5
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H87").Value = 4899.6
$ws.Range("I87").Value = 999.6667
$ws.Range("J87").Value = 10749.5
$ws.Range("K87").Value = 2999.0001
$ws.Range("L87").Value = 32248.5
$ws.Range("M87").Value = -1751.0001
$ws.Range("N87").Value = -34744.5
$ws.Range("H90").Value = 4899.6
$ws.Range("I90").Value = 999.6667
$ws.Range("J90").Value = 10749.5
$ws.Range("K90").Value = 8997.0003
$ws.Range("L90").Value = 96745.5
$ws.Range("M90").Value = -2757.0003
$ws.Range("N90").Value = -109225.5
$ws.Range("H94").Value = 7590.6
$ws.Range("I94").Value = 650
$ws.Range("J94").Value = 35353
$ws.Range("K94").Value = 1950
$ws.Range("L94").Value = 106059
$ws.Range("M94").Value = -1274
$ws.Range("N94").Value = -107411
$ws.Range("H108").Value = 355.375
$ws.Range("I108").Value = 355.375
$ws.Range("K108").Value = 1066.125
$ws.Range("M108").Value = 1813.875
$ws.Range("H116").Value = 14601923
$ws.Range("J116").Value = 257497.5
$ws.Range("L116").Value = 772492.5
$ws.Range("N116").Value = -779376.5
$ws.Range("H118").Value = 1763.3334
$ws.Range("I118").Value = 1750
$ws.Range("K118").Value = 5250
$ws.Range("M118").Value = -4007
$ws.Range("H119").Value = 21212
$ws.Range("J119").Value = 21212
$ws.Range("L119").Value = 63636
$ws.Range("N119").Value = -73312
$ws.Range("H122").Value = 6553.1665
$ws.Range("J122").Value = 9377.5
$ws.Range("L122").Value = 84397.5
$ws.Range("N122").Value = -89297.5
$ws.Range("H129").Value = 2225.2307
$ws.Range("J129").Value = 4224.8335
$ws.Range("L129").Value = 12674.5005
$ws.Range("N129").Value = -22674.5005
$ws.Range("H131").Value = 1371.495
$ws.Range("J131").Value = 1488.0581
$ws.Range("L131").Value = 4464.1743
$ws.Range("N131").Value = -14544.1743
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3093.1667
$ws.Range("J80").Value = 3499.75
$ws.Range("L80").Value = 3499.75
$ws.Range("N80").Value = -5495.75
$ws.Range("H83").Value = 3093.1667
$ws.Range("J83").Value = 3499.75
$ws.Range("L83").Value = 17498.75
$ws.Range("N83").Value = -27482.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2161.9412
$ws.Range("J16").Value = 2518
$ws.Range("L16").Value = 2518
$ws.Range("N16").Value = -2858
$ws.Range("H100").Value = 3665.4443
$ws.Range("I100").Value = 3298.3333
$ws.Range("J100").Value = 4399.6665
$ws.Range("K100").Value = 3298.3333
$ws.Range("L100").Value = 4399.6665
$ws.Range("M100").Value = -2757.3333
$ws.Range("N100").Value = -5481.6665
$ws.Range("H122").Value = 8464.579
$ws.Range("I122").Value = 8886.556
$ws.Range("J122").Value = 8084.8
$ws.Range("K122").Value = 26659.668
$ws.Range("L122").Value = 24254.4
$ws.Range("M122").Value = -24209.668
$ws.Range("N122").Value = -29154.4
$ws.Range("H132").Value = 1305854.6
$ws.Range("I132").Value = 5333.3335
$ws.Range("K132").Value = 16000.0005
$ws.Range("M132").Value = -13470.0005
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3920.3684
$ws.Range("I122").Value = 2967.9375
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 8903.8125
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -6453.8125
$ws.Range("N122").Value = -31900
$ws.Range("H136").Value = 487169.3
$ws.Range("I136").Value = 7600.3335
$ws.Range("K136").Value = 22801.0005
$ws.Range("M136").Value = -20251.0005
